$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 updates
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.25
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.73

# Row 9 updates
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.7

# Row 14 updates
$ws.Range("G14").Value = 2.47
$ws.Range("I14").Value = 2.62
$ws.Range("J14").Value = 3.05
$ws.Range("K14").Value = 2.12
$ws.Range("L14").Value = 3.2
$ws.Range("O14").Value = 1.34
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 2.02
$ws.Range("T14").Value = 2.75
$ws.Range("W14").Value = 7.9
$ws.Range("X14").Value = 12
$ws.Range("Y14").Value = 9.5
$ws.Range("Z14").Value = 26
$ws.Range("AA14").Value = 21
$ws.Range("AB14").Value = 32
$ws.Range("AH14").Value = 8.25
$ws.Range("AI14").Value = 13
$ws.Range("AJ14").Value = 10
$ws.Range("AK14").Value = 29
$ws.Range("AL14").Value = 23
$ws.Range("AM14").Value = 32
$ws.Range("AN14").Value = 4.45
$ws.Range("AO14").Value = 13
$ws.Range("AP14").Value = 21
$ws.Range("AQ14").Value = 55
$ws.Range("AR14").Value = 90
$ws.Range("AS14").Value = 250
$ws.Range("AT14").Value = 2.75
$ws.Range("AW14").Value = 4.6
$ws.Range("AX14").Value = 14
$ws.Range("AY14").Value = 22
$ws.Range("AZ14").Value = 60
$ws.Range("BA14").Value = 100
$ws.Range("BB14").Value = 300

$wb.Save()
